$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 10 new coded-segment rows (379-388) below the existing data (A1:M378).
#
# Columns A-D, H and L are identical for every new row to what row 378 already
# has (A="(bullet)", B/C blank, D="18833", H=0, L="Sonia"), so the cleanest way
# to get those (plus the shared cell formatting: fill/border/number-format/
# alignment) onto the new rows is to copy row 378's formats onto the new block
# first. Only E, F, G, I, J, K and M then need explicit per-row values.
#
# Document name (column D, "18833") and a few Segment (column I) values are
# bare digits ("2004"/"2006") which, if assigned via a normal .Value write,
# get auto-coerced to numbers instead of staying text. Force those specific
# cells to text first (while they still have the sheet's default formatting),
# then the row-378 format copy below restores the correct shared cell style
# on top without disturbing the text value/type.
for ($r = 379; $r -le 388; $r++) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = "18833"
}
$ws.Cells.Item(383, 9).NumberFormat = "@"
$ws.Cells.Item(384, 9).NumberFormat = "@"
$ws.Cells.Item(385, 9).NumberFormat = "@"
$ws.Cells.Item(388, 9).NumberFormat = "@"
$ws.Cells.Item(383, 9).Value = "2004"
$ws.Cells.Item(384, 9).Value = "2006"
$ws.Cells.Item(385, 9).Value = "2006"
$ws.Cells.Item(388, 9).Value = "2004"

$ws.Range("A378:M378").Copy()
$ws.Range("A379:M388").PasteSpecial(-4122)

# Columns A, D, H and L are the same on every new row as they are on row 378
# ("(bullet)", "18833", 0, "Sonia"); the format-only paste above does not carry
# values, so set them explicitly now.
for ($r = 379; $r -le 388; $r++) {
    $ws.Cells.Item($r, 1).Value = "●"
    $ws.Cells.Item($r, 4).Value = "18833"
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 12).Value = "Sonia"
}

$ws.Cells.Item(379, 5).Value = "Bacteria:Binomial (genus species)"
$ws.Cells.Item(379, 6).Value = "1: 709"
$ws.Cells.Item(379, 7).Value = "1: 722"
$ws.Cells.Item(379, 9).Value = ".  `npneumoniae"
$ws.Cells.Item(379, 10).Value = 13
$ws.Cells.Item(379, 11).Value = 0.040218
$ws.Cells.Item(379, 13).Value = "11/8/18 14:41:00"
$ws.Rows.Item(379).RowHeight = 30

$ws.Cells.Item(380, 5).Value = "Bacteria:Binomial (genus species)"
$ws.Cells.Item(380, 6).Value = "1: 727"
$ws.Cells.Item(380, 7).Value = "1: 740"
$ws.Cells.Item(380, 9).Value = " H. influenzae"
$ws.Cells.Item(380, 10).Value = 13
$ws.Cells.Item(380, 11).Value = 0.040218
$ws.Cells.Item(380, 13).Value = "11/8/18 14:41:00"
$ws.Rows.Item(380).RowHeight = 16

$ws.Cells.Item(381, 5).Value = "Event month"
$ws.Cells.Item(381, 6).Value = "2: 2174"
$ws.Cells.Item(381, 7).Value = "2: 2180"
$ws.Cells.Item(381, 9).Value = "October"
$ws.Cells.Item(381, 10).Value = 7
$ws.Cells.Item(381, 11).Value = 0.021656
$ws.Cells.Item(381, 13).Value = "11/13/18 08:38:00"
$ws.Rows.Item(381).RowHeight = 16

$ws.Cells.Item(382, 5).Value = "Event month"
$ws.Cells.Item(382, 6).Value = "2: 2192"
$ws.Cells.Item(382, 7).Value = "2: 2196"
$ws.Cells.Item(382, 9).Value = "March"
$ws.Cells.Item(382, 10).Value = 5
$ws.Cells.Item(382, 11).Value = 0.015468
$ws.Cells.Item(382, 13).Value = "11/13/18 08:38:00"
$ws.Rows.Item(382).RowHeight = 16

$ws.Cells.Item(383, 5).Value = "Event year"
$ws.Cells.Item(383, 6).Value = "2: 2182"
$ws.Cells.Item(383, 7).Value = "2: 2185"
$ws.Cells.Item(383, 10).Value = 4
$ws.Cells.Item(383, 11).Value = 0.012375
$ws.Cells.Item(383, 13).Value = "11/13/18 08:38:00"
$ws.Rows.Item(383).RowHeight = 16

$ws.Cells.Item(384, 5).Value = "Event year"
$ws.Cells.Item(384, 6).Value = "2: 2198"
$ws.Cells.Item(384, 7).Value = "2: 2201"
$ws.Cells.Item(384, 10).Value = 4
$ws.Cells.Item(384, 11).Value = 0.012375
$ws.Cells.Item(384, 13).Value = "11/13/18 08:38:00"
$ws.Rows.Item(384).RowHeight = 16

$ws.Cells.Item(385, 5).Value = "B"
$ws.Cells.Item(385, 6).Value = "2: 2198"
$ws.Cells.Item(385, 7).Value = "2: 2201"
$ws.Cells.Item(385, 10).Value = 4
$ws.Cells.Item(385, 11).Value = 0.012375
$ws.Cells.Item(385, 13).Value = "11/13/18 08:38:00"
$ws.Rows.Item(385).RowHeight = 16

$ws.Cells.Item(386, 5).Value = "B"
$ws.Cells.Item(386, 6).Value = "2: 2192"
$ws.Cells.Item(386, 7).Value = "2: 2196"
$ws.Cells.Item(386, 9).Value = "March"
$ws.Cells.Item(386, 10).Value = 5
$ws.Cells.Item(386, 11).Value = 0.015468
$ws.Cells.Item(386, 13).Value = "11/13/18 08:38:00"
$ws.Rows.Item(386).RowHeight = 16

$ws.Cells.Item(387, 5).Value = "B"
$ws.Cells.Item(387, 6).Value = "2: 2174"
$ws.Cells.Item(387, 7).Value = "2: 2180"
$ws.Cells.Item(387, 9).Value = "October"
$ws.Cells.Item(387, 10).Value = 7
$ws.Cells.Item(387, 11).Value = 0.021656
$ws.Cells.Item(387, 13).Value = "11/13/18 08:38:00"
$ws.Rows.Item(387).RowHeight = 16

$ws.Cells.Item(388, 5).Value = "B"
$ws.Cells.Item(388, 6).Value = "2: 2182"
$ws.Cells.Item(388, 7).Value = "2: 2185"
$ws.Cells.Item(388, 10).Value = 4
$ws.Cells.Item(388, 11).Value = 0.012375
$ws.Cells.Item(388, 13).Value = "11/13/18 08:38:00"
$ws.Rows.Item(388).RowHeight = 16
